$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.767.48"
$ws.Range("E2").Value = "  -3.90%  "
$ws.Range("D3").Value = "2.904.43"
$ws.Range("E3").Value = "  -4.23%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'589.82"
$ws.Range("E5").Value = "  -0.84%  "
$ws.Range("D6").Value = "'144.49"
$ws.Range("E6").Value = "  -6.00%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'0.504"
$ws.Range("E8").Value = "  -2.00%  "
$ws.Range("D9").Value = "2.904.20"
$ws.Range("E9").Value = "  -4.16%  "
$ws.Range("D10").Value = "'6.71"
$ws.Range("E10").Value = "  -4.89%  "
$ws.Range("E11").Value = "  -4.56%  "
$ws.Range("E12").Value = "  -4.47%  "
$ws.Range("E13").Value = "  -3.84%  "
$ws.Range("D14").Value = "'33.36"
$ws.Range("E14").Value = "  -6.75%  "
$ws.Range("D15").Value = "'0.126"
$ws.Range("E15").Value = "  +1.27%  "
$ws.Range("D16").Value = "3.386.76"
$ws.Range("D17").Value = "60.727.16"
$ws.Range("E17").Value = "  -3.93%  "
$ws.Range("E18").Value = "  -5.39%  "
$ws.Range("D19").Value = "2.904.83"
$ws.Range("E19").Value = "  -4.29%  "
$ws.Range("D20").Value = "'429.06"
$ws.Range("E20").Value = "  -4.82%  "
$ws.Range("E21").Value = "  -5.18%  "
$ws.Range("D22").Value = "'0.682"
$ws.Range("E22").Value = "  -2.06%  "
$ws.Range("D23").Value = "'7.06"
$ws.Range("E23").Value = "  -6.17%  "
$ws.Range("D24").Value = "'81.86"
$ws.Range("E24").Value = "  -1.72%  "
$ws.Range("D25").Value = "'10.81"
$ws.Range("E25").Value = "  -6.21%  "
$ws.Range("E26").Value = "  -5.36%  "
$ws.Range("D27").Value = "'11.94"
$ws.Range("E27").Value = "  -3.84%  "
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("D29").Value = "'2.23"
$ws.Range("E29").Value = "  -2.24%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("D32").Value = "'7.02"
$ws.Range("E32").Value = "  -7.36%  "
$ws.Range("D33").Value = "'26.47"
$ws.Range("E33").Value = "  -4.41%  "
$ws.Range("E34").Value = "  -4.10%  "
$ws.Range("D35").Value = "0.0₃0847"
$ws.Range("E35").Value = "  -3.03%  "
$ws.Range("E36").Value = "  -3.47%  "
$ws.Range("E37").Value = "  -5.39%  "
$ws.Range("D38").Value = "'3.00"
$ws.Range("E38").Value = "  -5.25%  "
$ws.Range("D39").Value = "'49.39"
$ws.Range("E39").Value = "  -2.41%  "
$ws.Range("D40").Value = "'0.125"
$ws.Range("E40").Value = "  -4.69%  "
$ws.Range("E41").Value = "  -5.76%  "
$ws.Range("D42").Value = "'8.59"
$ws.Range("E42").Value = "  -5.53%  "
$ws.Range("D43").Value = "'0.292"
$ws.Range("E43").Value = "  -5.46%  "
$ws.Range("D44").Value = "'40.63"
$ws.Range("E44").Value = "  -7.64%  "
$ws.Range("E45").Value = "  -3.19%  "
$ws.Range("D46").Value = "'373.68"
$ws.Range("E46").Value = "  -5.32%  "
$ws.Range("D47").Value = "2.698.72"
$ws.Range("E47").Value = "  -1.03%  "
$ws.Range("D48").Value = "'131.94"
$ws.Range("E48").Value = "  -1.30%  "
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("D50").Value = "'24.04"
$ws.Range("E50").Value = "  -9.00%  "
$ws.Range("E51").Value = "  -2.74%  "
